$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 41: 生活费 expense amount 300 -> 400 ---
$ws.Range("D41").Value = 400

# --- Add new row 42: 2018-06-10 支出 生活费 300, 生活费(6/10-6/19) ---
# Copy the formatting from row 41 first so the new row matches the
# existing table's styles (fills/borders/number formats), then overwrite
# the values.
$ws.Range("B41:G41").Copy()
$ws.Range("B42:G42").PasteSpecial(-4122)

$ws.Range("B42").Value = 40
$ws.Range("C42").Value = "支出"
$ws.Range("D42").Value = 300
$ws.Range("E42").Value = 43261
$ws.Range("F42").Value = "生活费"
$ws.Range("G42").Value = "生活费(6/10-6/19)"

# --- Update the view's current selection to K37 (matches the source edit) ---
$ws.Range("K37").Select()

# Recalculate so the summary formulas (I3/J3/K3/J9/...) pick up the new data
$excel.CalculateFull()
